$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6952002048492432
$ws.Range("B1").Value = 1.514212489128113
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.002576112747192
$ws.Range("E1").Value = 1.257044434547424
